# Case_4_180 (380 kV case): refresh res_bus/vm_pu.xlsx power-flow results.
# Bus 0 (column B) is re-solved with a 1.02 pu voltage setpoint (was 1.05 pu),
# so every downstream bus voltage in columns C:F and I:N is updated to the
# new load-flow solution for rows 2-25 (buses 0-23). Column G (always 1) is
# untouched, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Worksheet columns that carry per-bus voltage results, in sheet order.
$colLetters = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M", "N")
$colNumbers = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

# Each inner array is: row number, followed by the new value for each of the
# columns above, in the same order.
$rowData = @(
    @(2, 1.02, 1.036334033533864, 1.03863106661284, 1.050033362760966, 1.057809145830273, 1.038928330697977, 1.041442918201266, 1.041418550297654, 1.052788765023705, 1.060543099951803, 1.017755959234141),
    @(3, 1.02, 1.037148885291334, 1.039222024174919, 1.050839500471663, 1.058700438476468, 1.039105434834454, 1.041902369288703, 1.04182013963994, 1.053407284170311, 1.061248105452836, 1.017909832629264),
    @(4, 1.02, 1.037676836105266, 1.039604957942501, 1.051362188868775, 1.059278383898431, 1.039219202570195, 1.042199675442888, 1.042079840125888, 1.053807915018418, 1.061704871101677, 1.018009363364923),
    @(5, 1.02, 1.03789894944938, 1.039766071945229, 1.051582179737425, 1.059521642128399, 1.039266831002108, 1.042324664227439, 1.042188979961262, 1.053976436246721, 1.061897032526103, 1.018051197088238),
    @(6, 1.02, 1.03793625275567, 1.03979313118799, 1.051619131933181, 1.059562503202558, 1.039274816304753, 1.042345650432769, 1.042207302741454, 1.054004737287767, 1.061929305267852, 1.018058220613274),
    @(7, 1.02, 1.037679803356991, 1.039607110252408, 1.051365127407121, 1.059281633189875, 1.039219839768503, 1.042201345545697, 1.042081298609532, 1.053810166431983, 1.061707438235271, 1.018009922385122),
    @(8, 1.02, 1.036609273177214, 1.038830669646201, 1.050305579442256, 1.058110108665586, 1.038988355402803, 1.041598188902388, 1.041554300588349, 1.052997711130069, 1.060781238760307, 1.017807968568105),
    @(9, 1.02, 1.034728201156323, 1.037466731596302, 1.048446749486801, 1.056055159826786, 1.038574121112113, 1.040535486386607, 1.040624532118932, 1.051569251894725, 1.059153673683926, 1.017451847765995),
    @(10, 1.02, 1.033477844526655, 1.036560401690827, 1.047213176459441, 1.054691650001937, 1.03829375168985, 1.039827184824818, 1.040003999210201, 1.050619184320068, 1.058071768397853, 1.017214292009874),
    @(11, 1.02, 1.032937323311209, 1.036168675654407, 1.046680389272509, 1.054102791402141, 1.038171358226204, 1.039520537825916, 1.039735154013636, 1.050208345245738, 1.057604057189818, 1.017111399963826),
    @(12, 1.02, 1.032736685324942, 1.03602328126233, 1.046482694387483, 1.053884298013594, 1.038125747566546, 1.039406644603453, 1.039635271590648, 1.05005582501328, 1.057430444473538, 1.017073177417809),
    @(13, 1.02, 1.032779716687396, 1.036054463841747, 1.046525091258022, 1.053931154924291, 1.038135537906378, 1.039431074652809, 1.039656697655477, 1.050088537308036, 1.057467679711146, 1.017081376450702),
    @(14, 1.02, 1.032920735735246, 1.036156655052159, 1.04666404353798, 1.054084725866607, 1.038167591055572, 1.03951112318306, 1.039726898129675, 1.050195736159226, 1.057589703931103, 1.01710824055223),
    @(15, 1.02, 1.033007640252857, 1.036219633072918, 1.046749683979199, 1.054179377181573, 1.038187320431631, 1.039560444972003, 1.039770148167252, 1.050261796047185, 1.057664902510058, 1.017124791914632),
    @(16, 1.02, 1.033513736068996, 1.036586414598464, 1.047248564570989, 1.054730763440293, 1.038301853708245, 1.039847537186017, 1.040021838518751, 1.050646462005564, 1.058102825029851, 1.01722112004898),
    @(17, 1.02, 1.033831436490975, 1.036816681213555, 1.047561864166889, 1.055077049919294, 1.038373432411392, 1.040027637436918, 1.04017967784993, 1.050887900224841, 1.058377727198678, 1.01728153675678),
    @(18, 1.02, 1.034016831668549, 1.036951061259334, 1.047744737423374, 1.055279182395348, 1.038415087330198, 1.040132691822104, 1.040271728271274, 1.051028779671404, 1.058538146284358, 1.017316773986278),
    @(19, 1.02, 1.034080061150056, 1.03699689310288, 1.047807114621628, 1.055348129630816, 1.038429274322487, 1.040168513454391, 1.040303112573634, 1.051076824798493, 1.058592857427536, 1.017328788479215),
    @(20, 1.02, 1.033797341345239, 1.036791968590981, 1.047528236531238, 1.055039881186335, 1.038365762590431, 1.040008313862634, 1.040162744683637, 1.050861990734835, 1.058348225203197, 1.017275054903461),
    @(21, 1.02, 1.032879205381687, 1.036126559233949, 1.046623119852492, 1.054039496550724, 1.038158156289541, 1.039487550623875, 1.039706226419023, 1.050164166446292, 1.057553767622207, 1.017100329846623),
    @(22, 1.02, 1.032302722604943, 1.03570882785571, 1.046055229888532, 1.05341187516882, 1.038026768430753, 1.039160179541201, 1.039419072376574, 1.049725901627039, 1.057054932903765, 1.016990451176634),
    @(23, 1.02, 1.032608251974719, 1.035930214094924, 1.046356165382509, 1.053744459537061, 1.038096500668662, 1.039333719654926, 1.03957130940145, 1.049958187594534, 1.057319310321979, 1.017048701874612),
    @(24, 1.02, 1.033812747208799, 1.036803134947108, 1.047543431008092, 1.055056675675391, 1.038369228548252, 1.04001704533893, 1.040170396097846, 1.050873697955731, 1.058361555658788, 1.017277983782724),
    @(25, 1.02, 1.035213860196015, 1.03781882785922, 1.048926314641817, 1.056585284066963, 1.038681956313836, 1.040535486386607, 1.032819164395039, 1.051569251894725, 1.038461113387647, 1.013371780325615)
)

foreach ($entry in $rowData) {
    $rowNum = $entry[0]
    for ($i = 0; $i -lt $colNumbers.Count; $i++) {
        $colNum = $colNumbers[$i]
        $newValue = $entry[$i + 1]
        $ws.Cells.Item($rowNum, $colNum).Value = $newValue
    }
}

